# "fix hettinger ND, update query counties, muni"
#
# The source workbook is a lookup table of school-district records
# (columns: state.abb | name | id | ncesID | name_nces). This change
# removes a handful of stale/duplicate rows' data (clearing the row's
# content but leaving the row in place) and moves the sheet's active
# selection to reflect where the author was working afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 138 (ND / 108426 / 3800406) -> blank out the data, keep formatting
$ws.Range("A138").ClearContents()
$ws.Range("C138").ClearContents()
$ws.Range("D138").ClearContents()

# --- Row 163 (NJ / 1270501 / 3409990) -> blank out the data
# A163 loses its style entirely (full clear); C163/D163 keep formatting
$ws.Range("A163").Clear()
$ws.Range("C163").ClearContents()
$ws.Range("D163").ClearContents()

# --- Row 173 (SC / bamberg county school district / 1268144 / 4503916)
# duplicate of "Bamberg School District No. 1" already on row 137 -> clear it.
# B173 loses its style entirely (full clear); A173/C173/D173 keep formatting
$ws.Range("A173").ClearContents()
$ws.Range("B173").Clear()
$ws.Range("C173").ClearContents()
$ws.Range("D173").ClearContents()

# --- Row 178 (NJ / 1240267 / 3409990) -> blank out the data
$ws.Range("A178").ClearContents()
$ws.Range("C178").ClearContents()
$ws.Range("D178").ClearContents()

# --- Row 203 (orphan id 32029) -> blank out the data
$ws.Range("C203").ClearContents()

# --- Update the sheet's active selection / scroll position to where the
# author ended up after this pass (freeze-pane top row 167 -> 170, cursor
# moved from G191 down to B205).
$win = $excel.ActiveWindow
$win.ScrollRow = 170
$ws.Range("B205").Select() | Out-Null
